# Insert a new weekly price-observation record at row 23 (pushing the
# existing rows 23..92 down to 24..93), matching the diff which adds a
# new "Arveja Verde" record for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23; Excel shifts the old row 23
# (and everything below it) down by one, and the new row inherits the
# formatting (incl. the date style on column D) from the row above it.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record's data.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44600
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112022
$ws.Range("G23").Value = "Arveja Verde"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 22000
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región de La Araucanía"
$ws.Range("P23").Value = 880
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
